$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 0.55000000000000004
$ws.Range("K2").Value = 0.55000000000000004
$ws.Range("L2").Value = 0.75
$ws.Range("M2").Value = 0.75
$ws.Range("N2").Value = 1.25
$ws.Range("O2").Value = 1.85
$ws.Range("P2").Value = 1.85
$ws.Range("Q2").Value = 1.85
$ws.Range("R2").Value = 1.85

# Row 3
$ws.Range("D3").Value = 0.25
$ws.Range("E3").Value = 0.25
$ws.Range("F3").Value = 0.25
$ws.Range("G3").Value = 0.25
$ws.Range("H3").Value = 0.25
$ws.Range("I3").Value = 0.25
$ws.Range("J3").Value = 0.25
$ws.Range("K3").Value = 0.25
$ws.Range("L3").Value = 0.25
$ws.Range("M3").Value = 0.25
$ws.Range("N3").Value = 0.25
$ws.Range("O3").Value = 0.25
$ws.Range("P3").Value = 0.25
$ws.Range("Q3").Value = 0.25
$ws.Range("R3").Value = 0.5

# Row 4
$ws.Range("R4").Value = 2.875

# Row 5
$ws.Range("R5").Value = 4.375

# Row 7
$ws.Range("D7").Value = 0.35
$ws.Range("E7").Value = 0.35
$ws.Range("F7").Value = 0.85
$ws.Range("G7").Value = 0.85
$ws.Range("H7").Value = 0.85
$ws.Range("I7").Value = 0.85
$ws.Range("J7").Value = 0.85
$ws.Range("K7").Value = 0.95
$ws.Range("L7").Value = 0.95
$ws.Range("M7").Value = 0.95
$ws.Range("N7").Value = 0.95
$ws.Range("O7").Value = 0.95
$ws.Range("P7").Value = 0.95
$ws.Range("Q7").Value = 1.65
$ws.Range("R7").Value = 2

# Row 10
$ws.Range("R10").Value = 0.65

# Row 11
$ws.Range("J11").Value = 1.2999999999999998
$ws.Range("K11").Value = 1.2999999999999998
$ws.Range("L11").Value = 1.4999999999999998
$ws.Range("M11").Value = 1.4999999999999998
$ws.Range("N11").Value = 1.4999999999999998
$ws.Range("O11").Value = 1.4999999999999998
$ws.Range("P11").Value = 1.4999999999999998
$ws.Range("Q11").Value = 1.4999999999999998
$ws.Range("R11").Value = 1.4999999999999998

# Row 12
$ws.Range("R12").Value = 2.7499999999999996

# Update the active selection in the bottom-right frozen pane to C13
$ws.Range("C13").Select()
